$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.467323064804077
$ws.Range("B1").Value = 1.715182065963745
$ws.Range("C1").Value = 2.641993284225464
$ws.Range("D1").Value = 5.037627220153809
$ws.Range("E1").Value = 1.517630815505981
